$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: EFT - Şube ---
$ws.Range("F3").ClearContents()
$ws.Range("H3").ClearContents()

# --- Row 4: HESAPTAN EFT - ATM ---
$ws.Range("F4").ClearContents()
$ws.Range("H4").ClearContents()

# --- Row 5: HESAPTAN EFT - Mobil ---
$ws.Range("F5").ClearContents()
$ws.Range("H5").ClearContents()

# --- Row 6: DÜZENLİ EFT ---
$ws.Range("H6").ClearContents()

# --- Row 8: HESAPTAN HAVALE - Şube ---
$ws.Range("F8").ClearContents()
$ws.Range("H8").ClearContents()

# --- Row 9: HESAPTAN HAVALE - ATM ---
$ws.Range("F9").ClearContents()
$ws.Range("H9").ClearContents()

# --- Row 10: HESAPTAN HAVALE - Mobil ---
$ws.Range("F10").ClearContents()
$ws.Range("H10").ClearContents()

# --- Row 11: DÜZENLİ HAVALE ---
$ws.Range("H11").ClearContents()

# --- Row 13: GELEN SWIFT ---
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("F13").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# --- Row 14: GİDEN SWIFT - Mobil ---
$ws.Range("F14").ClearContents()
$ws.Range("H14").ClearContents()

# --- Row 15: ÇEK TAHSİLİ BAŞKA BANKA ---
$ws.Range("D15").Value = " Asgari Tutar:  Azami Tutar: "

# --- Row 17: AYNI ŞUBE ÇEK TAHSİLATI ---
$ws.Range("D17").Value = " Asgari Tutar:  Azami Tutar: "

# --- Row 20: ÇEK İADE ---
$ws.Range("D20").ClearContents()

# --- Row 21: BLOKE ÇEK DÜZENLEME ---
$ws.Range("D21").Value = " Asgari Tutar:  Azami Tutar: "

# --- Row 22: YP ÇEK TAKASA GÖNDERME ---
$ws.Range("D22").Value = " Asgari Tutar:  Azami Tutar: "

# --- Row 23: ÇEK KARNESİ SAYFA ÜCRETİ ---
$ws.Range("D23").ClearContents()

# --- Row 24: SENET TAHSİLE ALMA ---
$ws.Range("D24").ClearContents()

# --- Row 25: MUAMELESİZ SENET İADESİ ---
$ws.Range("D25").ClearContents()
